$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.240.10"
$ws.Range("E2").Value = '  +5.51%  '
$ws.Range("D3").Value = "'1.916.03"
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  -0.62%  '
$ws.Range("D5").Value = "'329.57"
$ws.Range("E5").Value = '  +4.74%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("D7").Value = "'0.5196"
$ws.Range("E7").Value = '  +2.28%  '
$ws.Range("D8").Value = "'0.4065"
$ws.Range("E8").Value = '  +3.85%  '
$ws.Range("D9").Value = "'0.08512"
$ws.Range("E9").Value = '  +1.60%  '
$ws.Range("D10").Value = "'1.128"
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("D11").Value = "'42.80"
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").Value = "'23.37"
$ws.Range("E12").Value = '  +14.72%  '
$ws.Range("D14").Value = "'1.917.07"
$ws.Range("E14").Value = '  +2.33%  '
$ws.Range("D15").Value = "'7.400"
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").Value = "'95.28"
$ws.Range("E17").Value = '  +2.29%  '
$ws.Range("D18").Value = "'0.00001114"
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("D19").Value = "'0.06700"
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").Value = "'18.54"
$ws.Range("E20").Value = '  +5.08%  '
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").Value = "'6.024"
$ws.Range("E22").Value = '  +1.56%  '
$ws.Range("D23").Value = "'30.243.14"
$ws.Range("E23").Value = '  +5.43%  '
$ws.Range("D24").Value = "'11.37"
$ws.Range("E24").Value = '  +2.54%  '
$ws.Range("D25").Value = "'2.228"
$ws.Range("E25").Value = '  +1.61%  '
$ws.Range("D26").Value = "'2.136.57"
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("D27").Value = "'21.51"
$ws.Range("E27").Value = '  +4.35%  '
$ws.Range("D28").Value = "'162.57"
$ws.Range("E28").Value = '  +3.24%  '
$ws.Range("D29").Value = "'2.411"
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("D30").Value = "'129.13"
$ws.Range("E30").Value = '  +2.03%  '
$ws.Range("D31").Value = "'1.111"
$ws.Range("E31").Value = '  +6.11%  '
$ws.Range("D32").Value = "'0.1068"
$ws.Range("E32").Value = '  +2.86%  '
$ws.Range("D33").Value = "'6.024"
$ws.Range("E33").Value = '  +4.22%  '
$ws.Range("D34").Value = "'3.658"
$ws.Range("E34").Value = '  +0.73%  '
$ws.Range("D35").Value = "'0.02496"
$ws.Range("E35").Value = '  +1.68%  '
$ws.Range("D36").Value = "'0.06579"
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").Value = "'0.2215"
$ws.Range("E37").Value = '  +2.38%  '
$ws.Range("D38").Value = "'5.207"
$ws.Range("E38").Value = '  +3.26%  '
$ws.Range("D39").Value = "'1.232"
$ws.Range("E39").Value = '  +3.33%  '
$ws.Range("D40").Value = "'12.02"
$ws.Range("E40").Value = '  +7.98%  '
$ws.Range("D41").Value = "'8.827"
$ws.Range("E41").Value = '  -2.10%  '
$ws.Range("D42").Value = "'0.6539"
$ws.Range("E42").Value = '  +2.35%  '
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("D44").Value = "'0.6150"
$ws.Range("E44").Value = '  +2.62%  '
$ws.Range("D45").Value = "'13.27"
$ws.Range("E45").Value = '  +1.75%  '
$ws.Range("D46").Value = "'3.740"
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("D47").Value = "'2.080"
$ws.Range("E47").Value = '  +3.76%  '
$ws.Range("D48").Value = "'1.248"
$ws.Range("E48").Value = '  +2.18%  '
$ws.Range("D49").Value = "'124.66"
$ws.Range("E49").Value = '  +2.03%  '
$ws.Range("E50").Value = '  -4.24%  '
$ws.Range("D51").Value = "'79.65"
$ws.Range("E51").Value = '  +4.65%  '
